$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("A").ColumnWidth = 24.29
